# Updated cryptos list values (Price column D, Volume(1h) column E).
# D-column values are plain decimal numbers in source data, but must
# be stored as literal text (matching the original inline-string cells).
# Temporarily force a text NumberFormat while assigning so Excel does not
# auto-parse e.g. "22.90" into the number 22.9, then restore the default
# "Normal" style so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.897.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.370.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.472"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  -5.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.947.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.375.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.983.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -5.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.189"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.40%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -6.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "167.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.407.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0760"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.99%  "
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.455.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.22%  "
$ws.Range("E49").Value = "  -4.99%  "
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("E51").Value = "  -3.89%  "
